# Convert the "Starting time" column (E2:E53) from a numeric Excel time
# serial value to plain text (e.g. 0.645833... -> "15:30"), matching the
# commit "Change class time to text".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$endRow = 53
$col = "E"

for ($row = $startRow; $row -le $endRow; $row++) {
    $cell = $ws.Range("$col$row")

    # Read the current (numeric) time-of-day value, e.g. 0.6458333333333333
    $fraction = [double]$cell.Value2

    $totalMinutes = [Math]::Round($fraction * 24 * 60)
    $hours = [Math]::Floor($totalMinutes / 60)
    $minutes = $totalMinutes % 60
    $timeText = "{0}:{1:d2}" -f $hours, $minutes

    # Switch the cell to text format and write the literal time string.
    $cell.NumberFormat = "@"
    $cell.Value = $timeText
}
